# Natmi following Dr Hou advice
# Rewrites the LR-pair result rows (2-10) of Sheet1 with updated values and
# extends the table from 6 data rows (the 3 sender/receiver combos that
# excluded "sCs" as a sender) to the full 3x3 grid of
# Sending cluster x Target cluster combinations (ECs / FAPs / sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A                B         C        D               E  F  G                    H           I                     J                     K  L                     M                    N           O                     P                     Q                    R                   S                     T
$data = @(
    @("ECs",  "Lrpap1", "Sorl1", "ECs",  3, 1, 6.540008666666668,  19.620026, 0.2365207520404831, 0.2365207520404831, 2, 0.6666666666666666, 86.09667233333334,  258.290017, 0.9236051571299395,  0.9236051571299394,  563.0729832311604,  5067.656849080444, 0.2184517863528418,  0.2184517863528418),
    @("ECs",  "Lrpap1", "Sorl1", "FAPs", 3, 1, 6.540008666666668,  19.620026, 0.2365207520404831, 0.2365207520404831, 3, 1,                   4.911922333333334,  14.735767,  0.05269282395635592, 0.05269282395635591, 32.12401462999357,  289.116131669942,  0.01246294634929408, 0.01246294634929408),
    @("ECs",  "Lrpap1", "Sorl1", "sCs",  3, 1, 6.540008666666668,  19.620026, 0.2365207520404831, 0.2365207520404831, 3, 1,                   2.209456,           6.628368,   0.02370201891370452, 0.02370201891370452, 14.44986138861867,  130.048752497568,  0.005606019338347147, 0.005606019338347147),
    @("FAPs", "Lrpap1", "Sorl1", "ECs",  3, 1, 14.67485766666667,  44.024573, 0.5307192311682535, 0.5307192311682536, 2, 0.6666666666666666, 86.09667233333334,  258.290017, 0.9236051571299395,  0.9236051571299394,  1263.456412065305,  11371.10770858774, 0.4901750188950355,  0.4901750188950355),
    @("FAPs", "Lrpap1", "Sorl1", "FAPs", 3, 1, 14.67485766666667,  44.024573, 0.5307192311682535, 0.5307192311682536, 3, 1,                   4.911922333333334,  14.735767,  0.05269282395635592, 0.05269282395635591, 72.0817611113879,   648.7358500024911, 0.02796509501820134, 0.02796509501820134),
    @("FAPs", "Lrpap1", "Sorl1", "sCs",  3, 1, 14.67485766666667,  44.024573, 0.5307192311682535, 0.5307192311682536, 3, 1,                   2.209456,           6.628368,   0.02370201891370452, 0.02370201891370452, 32.42345232076267,  291.811070886864,  0.01257911725501666, 0.01257911725501667),
    @("sCs",  "Lrpap1", "Sorl1", "ECs",  3, 1, 6.436020999999999,  19.308063, 0.2327600167912634, 0.2327600167912634, 2, 0.6666666666666666, 86.09667233333334,  258.290017, 0.9236051571299395,  0.9236051571299394,  554.1199911674523,  4987.079920507071, 0.2149783518820622,  0.2149783518820621),
    @("sCs",  "Lrpap1", "Sorl1", "FAPs", 3, 1, 6.436020999999999,  19.308063, 0.2327600167912634, 0.2327600167912634, 3, 1,                   4.911922333333334,  14.735767,  0.05269282395635592, 0.05269282395635591, 31.61323528770233,  284.519117589321,  0.01226478258886049, 0.01226478258886049),
    @("sCs",  "Lrpap1", "Sorl1", "sCs",  3, 1, 6.436020999999999,  19.308063, 0.2327600167912634, 0.2327600167912634, 3, 1,                   2.209456,           6.628368,   0.02370201891370452, 0.02370201891370452, 14.220105214576,    127.980946931184,  0.005516882320340706, 0.005516882320340706)
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowIndex = 2
foreach ($row in $data) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $rowIndex).Value = $row[$i]
    }
    $rowIndex++
}
